# DUA: Clarify data model, add simple design diagram
# Adds a new bulleted list item after "Use async throughout" describing
# the FK relationship between stock and trade.

$d = $word.ActiveDocument

# Locate the existing last bullet ("Use async throughout") in the
# "MVP Scaffold" list (numId = 1) and append a new list item after it,
# inheriting the same paragraph/numbering/run formatting.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.InsertAfter("FK Reference to stock ID in trade, one to many relationship between stock and trade")

Write-Output "Inserted new list item."
